$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet gets a new "latest reading" prepended above the existing
# date/price history every day, pushing the older rows down by one.
# Shift the existing data block (rows 2-24) down to rows 3-25, copying
# values/formats as-is so the historical rows are untouched other than
# their row number.
$ws.Range("A2:D24").Copy($ws.Range("A3:D25"))

# Fill in the new top row with the latest date and the same price
# values used throughout the rest of the table. Force column A to be
# stored as plain text (matching the rest of the date column) instead
# of letting Excel auto-convert the date-like string into a date
# serial number.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-12-14"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# Drop the (unformatted) new row back to the same plain styling used by
# the rest of the data rows.
$ws.Range("A2:D2").ClearFormats()
